$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Input")

# Columns A (발주일자) and B (납기일자) hold dates written as literal text
# (e.g. "2025-09-10"); force text format first so Excel does not auto-convert
# them into date serial numbers.
$ws.Range("A1:B7").NumberFormat = "@"

# --- Header row (row 1): new column labels ---
$ws.Cells.Item(1,1).Value = '발주일자'
$ws.Cells.Item(1,2).Value = '납기일자'
$ws.Cells.Item(1,3).Value = '거래처명'
$ws.Cells.Item(1,4).Value = '거래처 이메일'
$ws.Cells.Item(1,5).Value = '납품처명'
$ws.Cells.Item(1,6).Value = '납품처 이메일'
$ws.Cells.Item(1,7).Value = '프로젝트명'
$ws.Cells.Item(1,8).Value = '대분류'
$ws.Cells.Item(1,9).Value = '중분류'
$ws.Cells.Item(1,10).Value = '소분류'
$ws.Cells.Item(1,11).Value = '품목명'
$ws.Cells.Item(1,12).Value = '규격'
$ws.Cells.Item(1,13).Value = '수량'
$ws.Cells.Item(1,14).Value = '단가'
$ws.Cells.Item(1,15).Value = '총금액'
$ws.Cells.Item(1,16).Value = '비고'

# --- Data row 2 ---
$ws.Cells.Item(2,1).Value = '2025-09-10'
$ws.Cells.Item(2,2).Value = '2025-10-13'
$ws.Cells.Item(2,3).Value = '신호수'
$ws.Cells.Item(2,4).Value = '신호수@example.com'
$ws.Cells.Item(2,5).Value = '힐스테이트 도곡동1차'
$ws.Cells.Item(2,6).Value = 'delivery@example.com'
$ws.Cells.Item(2,7).Value = '힐스테이트 도곡동1차'
$ws.Cells.Item(2,8).Value = '4. 장비비'
$ws.Cells.Item(2,9).Value = '2) 신호수 외'
$ws.Cells.Item(2,10).Value = '기타'
$ws.Cells.Item(2,11).Value = '5월'
$ws.Cells.Item(2,12).Value = 'KS규격-1'
$ws.Cells.Item(2,13).Value = 7.2
$ws.Cells.Item(2,14).Value = 155000
$ws.Cells.Item(2,15).Value = 1227600
$ws.Cells.Item(2,16).Value = '서지원 '

# --- Data row 3 ---
$ws.Cells.Item(3,1).Value = '2025-09-16'
$ws.Cells.Item(3,2).Value = '2025-10-04'
$ws.Cells.Item(3,3).Value = '신호수'
$ws.Cells.Item(3,4).Value = '신호수@example.com'
$ws.Cells.Item(3,5).Value = '힐스테이트 도곡동1차'
$ws.Cells.Item(3,6).Value = 'delivery@example.com'
$ws.Cells.Item(3,7).Value = '힐스테이트 도곡동1차'
$ws.Cells.Item(3,8).Value = '4. 장비비'
$ws.Cells.Item(3,9).Value = '2) 신호수 외'
$ws.Cells.Item(3,10).Value = '기타'
$ws.Cells.Item(3,11).Value = '5월'
$ws.Cells.Item(3,12).Value = 'KS규격-2'
$ws.Cells.Item(3,13).Value = 6
$ws.Cells.Item(3,14).Value = 155000
$ws.Cells.Item(3,15).Value = 1023000
$ws.Cells.Item(3,16).Value = '탁영롱 '

# --- Data row 4 ---
$ws.Cells.Item(4,1).Value = '2025-08-31'
$ws.Cells.Item(4,2).Value = '2025-09-17'
$ws.Cells.Item(4,3).Value = '신호수'
$ws.Cells.Item(4,4).Value = '신호수@example.com'
$ws.Cells.Item(4,5).Value = '힐스테이트 도곡동1차'
$ws.Cells.Item(4,6).Value = 'delivery@example.com'
$ws.Cells.Item(4,7).Value = '힐스테이트 도곡동1차'
$ws.Cells.Item(4,8).Value = '4. 장비비'
$ws.Cells.Item(4,9).Value = '2) 신호수 외'
$ws.Cells.Item(4,10).Value = '기타'
$ws.Cells.Item(4,11).Value = '4월'
$ws.Cells.Item(4,12).Value = 'KS규격-3'
$ws.Cells.Item(4,13).Value = 2
$ws.Cells.Item(4,14).Value = 155000
$ws.Cells.Item(4,15).Value = 341000
$ws.Cells.Item(4,16).Value = '탁영롱'

# --- Data row 5 ---
$ws.Cells.Item(5,1).Value = '2025-09-17'
$ws.Cells.Item(5,2).Value = '2025-10-03'
$ws.Cells.Item(5,3).Value = '신호수'
$ws.Cells.Item(5,4).Value = '신호수@example.com'
$ws.Cells.Item(5,5).Value = '힐스테이트 도곡동1차'
$ws.Cells.Item(5,6).Value = 'delivery@example.com'
$ws.Cells.Item(5,7).Value = '힐스테이트 도곡동1차'
$ws.Cells.Item(5,8).Value = '4. 장비비'
$ws.Cells.Item(5,9).Value = '2) 신호수 외'
$ws.Cells.Item(5,10).Value = '기타'
$ws.Cells.Item(5,11).Value = '5월'
$ws.Cells.Item(5,12).Value = 'KS규격-4'
$ws.Cells.Item(5,13).Value = 1
$ws.Cells.Item(5,14).Value = 155000
$ws.Cells.Item(5,15).Value = 170500
$ws.Cells.Item(5,16).Value = '김병호 '

# --- Data row 6 ---
$ws.Cells.Item(6,1).Value = '2025-08-26'
$ws.Cells.Item(6,2).Value = '2025-09-17'
$ws.Cells.Item(6,3).Value = '신호수'
$ws.Cells.Item(6,4).Value = '신호수@example.com'
$ws.Cells.Item(6,5).Value = '힐스테이트 도곡동1차'
$ws.Cells.Item(6,6).Value = 'delivery@example.com'
$ws.Cells.Item(6,7).Value = '힐스테이트 도곡동1차'
$ws.Cells.Item(6,8).Value = '4. 장비비'
$ws.Cells.Item(6,9).Value = '2) 신호수 외'
$ws.Cells.Item(6,10).Value = '기타'
$ws.Cells.Item(6,11).Value = '5월'
$ws.Cells.Item(6,12).Value = 'KS규격-5'
$ws.Cells.Item(6,13).Value = 2
$ws.Cells.Item(6,14).Value = 155000
$ws.Cells.Item(6,15).Value = 341000
$ws.Cells.Item(6,16).Value = '정성식 '

# --- Data row 7 ---
$ws.Cells.Item(7,1).Value = '2025-09-12'
$ws.Cells.Item(7,2).Value = '2025-09-25'
$ws.Cells.Item(7,3).Value = '신호수'
$ws.Cells.Item(7,4).Value = '신호수@example.com'
$ws.Cells.Item(7,5).Value = '힐스테이트 도곡동1차'
$ws.Cells.Item(7,6).Value = 'delivery@example.com'
$ws.Cells.Item(7,7).Value = '힐스테이트 도곡동1차'
$ws.Cells.Item(7,8).Value = '4. 장비비'
$ws.Cells.Item(7,9).Value = '2) 신호수 외'
$ws.Cells.Item(7,10).Value = '기타'
$ws.Cells.Item(7,11).Value = '5월'
$ws.Cells.Item(7,12).Value = 'KS규격-6'
$ws.Cells.Item(7,13).Value = 1
$ws.Cells.Item(7,14).Value = 155000
$ws.Cells.Item(7,15).Value = 170500
$ws.Cells.Item(7,16).Value = '박수진 '

# Strip all cell styling from A1:P7 (the header used to be bold/bordered via
# style index 1 and the text-forced date cells picked up a quote-prefix style;
# the target layout has no styling on any of these cells) while preserving the
# underlying values and types.
$ws.Range("A1:P7").ClearFormats()

# Drop the now-unused column Q entirely (content + formatting) so the used
# range shrinks from A1:Q7 to A1:P7, matching the new 16-column layout.
$ws.Range("Q1:Q7").Clear()
